$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Duplicate the last weekly block (rows 102:110, the "2018.10.31" week)
# twice to create two new weekly blocks at rows 111:119 and 120:128,
# matching the same table layout (date band, column headers, 4 member
# rows, blank spacer row, "summary" band, blank row).
#
# NOTE: rows 102 and 109:110 are merged cells (A102:D102, A109:D110).
# Copying a multi-row range that includes a merged source cell causes
# this COM host to "split" the merged cell's single border style into
# several new per-edge border styles on paste (a correct visual result,
# but it bloats styles.xml with near-duplicate cellXfs entries instead
# of reusing the existing ones). Pre-merging the (still empty) target
# cells *before* pasting formats avoids that: the paste then lands on
# an already-merged destination and reuses the original style indices
# verbatim, matching how the author's saved file looks.
# ---------------------------------------------------------------------

$ws.Range("A111:D111").Merge()
$ws.Range("A118:D119").Merge()
$ws.Range("A120:D120").Merge()
$ws.Range("A127:D128").Merge()

$ws.Range("A102:D110").Copy()
$ws.Range("A111:D119").PasteSpecial(-4122)

$ws.Range("A102:D110").Copy()
$ws.Range("A120:D128").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Fill in the values for the new blocks. Dates are written before the
# "编写pc前端界面" task text so new shared-string entries land in the
# same order as the source workbook (66, 67, 68).
# ---------------------------------------------------------------------

# Block 1 / Block 2 header dates
$ws.Range("A111").Value = "日期：2018.11.5 第十周周一"
$ws.Range("A120").Value = "日期：2018.11.7 第十周周三"

# Block 1 column headers
$ws.Range("A112").Value = "组员"
$ws.Range("B112").Value = "计划内容"
$ws.Range("C112").Value = "完成情况"
$ws.Range("D112").Value = "备注"

# Block 2 column headers
$ws.Range("A121").Value = "组员"
$ws.Range("B121").Value = "计划内容"
$ws.Range("C121").Value = "完成情况"
$ws.Range("D121").Value = "备注"

# Block 1 member rows
$ws.Range("A113").Value = "余舒章"
$ws.Range("B113").Value = "编写pc前端界面"
$ws.Range("C113").Value = "进行中"

$ws.Range("A114").Value = "王嘉宇"
$ws.Range("B114").Value = "编写完善数据库"
$ws.Range("C114").Value = "进行中"

$ws.Range("A115").Value = "许俊杰"
$ws.Range("B115").Value = "使用springmvc、mybatis进行框架搭建"
$ws.Range("C115").Value = "进行中"

$ws.Range("A116").Value = "庞森杰"
$ws.Range("B116").Value = "使用安卓原生代码、百度地图sdk进行安卓端前端设计"
$ws.Range("C116").Value = "进行中"

# Block 1 summary band
$ws.Range("A118").Value = "总结："

# Block 2 member rows
$ws.Range("A122").Value = "余舒章"
$ws.Range("B122").Value = "编写pc前端界面"
$ws.Range("C122").Value = "进行中"

$ws.Range("A123").Value = "王嘉宇"
$ws.Range("B123").Value = "编写完善数据库"
$ws.Range("C123").Value = "进行中"

$ws.Range("A124").Value = "许俊杰"
$ws.Range("B124").Value = "使用springmvc、mybatis进行框架搭建"
$ws.Range("C124").Value = "进行中"

$ws.Range("A125").Value = "庞森杰"
$ws.Range("B125").Value = "使用安卓原生代码、百度地图sdk进行安卓端前端设计"
$ws.Range("C125").Value = "进行中"

# Block 2 summary band
$ws.Range("A127").Value = "总结："

# ---------------------------------------------------------------------
# Update window view state (scroll position + active selection) to
# match the author's saved view.
# ---------------------------------------------------------------------
$ws.Select()
$ws.Range("H107").Select()
